$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ34544253",
    "summ34676264",
    "summ34796864",
    "summ34918524",
    "summ35041970",
    "summ35155391",
    "summ35254283",
    "summ35464345",
    "summ35574141",
    "summ35671073",
    "summ35768850",
    "summ35866983",
    "summ35963809",
    "summ36110990",
    "summ36213037",
    "summ36313171",
    "summ36418363",
    "summ36518434",
    "summ36621810",
    "summ36722705",
    "summ36827412",
    "summ36936927",
    "summ37061703",
    "summ37195712",
    "summ37351705",
    "summ37493581",
    "summ37648861",
    "summ37797248",
    "summ37945716",
    "summ38146893",
    "summ38313659",
    "summ38448338",
    "summ38594932",
    "summ38721270",
    "summ38844837",
    "summ38981735",
    "summ39119011",
    "summ39282102",
    "summ39432761",
    "summ39574599",
    "summ39713386",
    "summ39845414",
    "summ39975554",
    "summ40107683",
    "summ40236447",
    "summ40371405",
    "summ40525236",
    "summ40663287",
    "summ40789886",
    "summ40919190"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
